# Scheduled-runner style refresh of market/profit figures (currentAveragePrice,
# currentAveragePriceNQ/HQ, LevePriceNQ/HQ, LeveProfitNQ/HQ columns H:N) across
# the ALC/ARM/BSM/CRP/CUL/GSM/LTW/WVR sheets, matching the upstream data pull.

$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H32").Value = 2983.8572
$ws.Range("I32").Value = 2631.3333
$ws.Range("J32").Value = 3248.25
$ws.Range("K32").Value = 2631.3333
$ws.Range("L32").Value = 3248.25
$ws.Range("M32").Value = -2305.3333
$ws.Range("N32").Value = -3900.25

$ws.Range("H62").Value = 0
$ws.Range("J62").Value = 0
$ws.Range("L62").Value = 0
$ws.Range("N62").Value = $null

$ws.Range("H65").Value = 0
$ws.Range("J65").Value = 0
$ws.Range("L65").Value = 0
$ws.Range("N65").Value = $null

$ws.Range("H103").Value = 734
$ws.Range("I103").Value = 1009.75
$ws.Range("J103").Value = 488.8889
$ws.Range("K103").Value = 3029.25
$ws.Range("L103").Value = 1466.6667
$ws.Range("M103").Value = -2443.25
$ws.Range("N103").Value = -2638.6667

$ws.Range("H106").Value = 3247.25
$ws.Range("I106").Value = 2663
$ws.Range("K106").Value = 2663
$ws.Range("M106").Value = -2032

$ws.Range("H116").Value = 5111.8887
$ws.Range("I116").Value = 4401
$ws.Range("K116").Value = 4401
$ws.Range("M116").Value = -959

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H36").Value = 2999
$ws.Range("I36").Value = 2999
$ws.Range("K36").Value = 2999
$ws.Range("M36").Value = -2653

$ws.Range("H74").Value = 100002640
$ws.Range("I74").Value = 100002640
$ws.Range("K74").Value = 100002640
$ws.Range("M74").Value = -100001766

$ws.Range("H77").Value = 100002640
$ws.Range("I77").Value = 100002640
$ws.Range("K77").Value = 500013200
$ws.Range("M77").Value = -500008832

$ws.Range("H97").Value = 935
$ws.Range("I97").Value = 935
$ws.Range("K97").Value = 935
$ws.Range("M97").Value = -439

$ws.Range("H102").Value = 6667311.5
$ws.Range("I102").Value = 7143491.5
$ws.Range("K102").Value = 7143491.5
$ws.Range("M102").Value = -7141869.5

$ws.Range("H110").Value = 48193
$ws.Range("I110").Value = 50577.7
$ws.Range("K110").Value = 50577.7
$ws.Range("M110").Value = -48532.7

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H20").Value = 2952.6924
$ws.Range("I20").Value = 2108.8
$ws.Range("K20").Value = 2108.8
$ws.Range("M20").Value = -1861.8

$ws.Range("H105").Value = 3114.8572
$ws.Range("I105").Value = 2258
$ws.Range("K105").Value = 2258
$ws.Range("M105").Value = -511

$ws.Range("H107").Value = 39032.406
$ws.Range("I107").Value = 2115.44
$ws.Range("K107").Value = 2115.44
$ws.Range("M107").Value = -195.4400000000001

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H16").Value = 1568079.1
$ws.Range("I16").Value = 1829175.6
$ws.Range("K16").Value = 1829175.6
$ws.Range("M16").Value = -1828888.6

$ws.Range("H22").Value = 625.4761999999999
$ws.Range("I22").Value = 625.4761999999999
$ws.Range("K22").Value = 625.4761999999999
$ws.Range("M22").Value = -275.4761999999999

$ws.Range("H31").Value = 3756.158
$ws.Range("I31").Value = 3663.5925
$ws.Range("K31").Value = 3663.5925
$ws.Range("M31").Value = -3368.5925

$ws.Range("H34").Value = 3756.158
$ws.Range("I34").Value = 3663.5925
$ws.Range("K34").Value = 3663.5925
$ws.Range("M34").Value = -3461.5925

$ws.Range("H92").Value = 89999
$ws.Range("J92").Value = 89999
$ws.Range("L92").Value = 89999
$ws.Range("N92").Value = -94991

$ws.Range("H113").Value = 1568079.1
$ws.Range("I113").Value = 1829175.6
$ws.Range("K113").Value = 1829175.6
$ws.Range("M113").Value = -1827005.6

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H2").Value = 266.15384
$ws.Range("J2").Value = 539.6667
$ws.Range("L2").Value = 3238.0002
$ws.Range("N2").Value = -3464.0002

$ws.Range("H80").Value = 850
$ws.Range("J80").Value = 850
$ws.Range("L80").Value = 2550
$ws.Range("N80").Value = -4422

$ws.Range("H83").Value = 850
$ws.Range("J83").Value = 850
$ws.Range("L83").Value = 7650
$ws.Range("N83").Value = -17010

$ws.Range("H122").Value = 684.4545000000001
$ws.Range("I122").Value = 489.66666
$ws.Range("K122").Value = 4406.99994
$ws.Range("M122").Value = -1956.99994

$ws.Range("H131").Value = 1109.8
$ws.Range("I131").Value = 786.8570999999999
$ws.Range("J131").Value = 1863.3334
$ws.Range("K131").Value = 2360.5713
$ws.Range("L131").Value = 5590.0002
$ws.Range("M131").Value = 2679.4287
$ws.Range("N131").Value = -15670.0002

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H80").Value = 2426.1667
$ws.Range("J80").Value = 2473.3333
$ws.Range("L80").Value = 2473.3333
$ws.Range("N80").Value = -4469.3333

$ws.Range("H83").Value = 2426.1667
$ws.Range("J83").Value = 2473.3333
$ws.Range("L83").Value = 12366.6665
$ws.Range("N83").Value = -22350.6665

$ws.Range("H113").Value = 58054
$ws.Range("I113").Value = 68911.734
$ws.Range("J113").Value = 3765.3333
$ws.Range("K113").Value = 68911.734
$ws.Range("L113").Value = 3765.3333
$ws.Range("M113").Value = -66741.734
$ws.Range("N113").Value = -8105.3333

$ws.Range("H126").Value = 2209.1765
$ws.Range("I126").Value = 2116
$ws.Range("K126").Value = 6348
$ws.Range("M126").Value = -3878

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value = 4531.5
$ws.Range("I7").Value = 4348.4
$ws.Range("J7").Value = 5447
$ws.Range("K7").Value = 4348.4
$ws.Range("L7").Value = 5447
$ws.Range("M7").Value = -4236.4
$ws.Range("N7").Value = -5671

$ws.Range("H16").Value = 3302.9
$ws.Range("I16").Value = 3012.8572
$ws.Range("K16").Value = 3012.8572
$ws.Range("M16").Value = -2842.8572

$ws.Range("H55").Value = 562.3333
$ws.Range("I55").Value = 599.5
$ws.Range("K55").Value = 599.5
$ws.Range("M55").Value = -426.5

$ws.Range("H93").Value = 3296.9167
$ws.Range("I93").Value = 2512.7778
$ws.Range("K93").Value = 2512.7778
$ws.Range("M93").Value = -1264.7778

$ws.Range("H100").Value = 27070786
$ws.Range("I100").Value = 31582444
$ws.Range("K100").Value = 31582444
$ws.Range("M100").Value = -31581903

$ws.Range("H126").Value = 4531.5
$ws.Range("I126").Value = 4348.4
$ws.Range("J126").Value = 5447
$ws.Range("K126").Value = 13045.2
$ws.Range("L126").Value = 16341
$ws.Range("M126").Value = -10575.2
$ws.Range("N126").Value = -21281

$ws.Range("H132").Value = 6417289.5
$ws.Range("I132").Value = 10006540
$ws.Range("K132").Value = 30019620
$ws.Range("M132").Value = -30017090

$ws.Range("H136").Value = 3199.1667
$ws.Range("I136").Value = 3287.4119
$ws.Range("K136").Value = 9862.235700000001
$ws.Range("M136").Value = -7312.235700000001

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H100").Value = 815.3333
$ws.Range("I100").Value = 1245.3334
$ws.Range("J100").Value = 385.33334
$ws.Range("K100").Value = 2490.6668
$ws.Range("L100").Value = 770.66668
$ws.Range("M100").Value = -1949.6668
$ws.Range("N100").Value = -1852.66668

$ws.Range("H104").Value = 26182.5
$ws.Range("J104").Value = 26182.5
$ws.Range("L104").Value = 26182.5
$ws.Range("N104").Value = -33170.5

$ws.Range("H122").Value = 1639.7
$ws.Range("I122").Value = 1556.1875
$ws.Range("K122").Value = 4668.5625
$ws.Range("M122").Value = -2218.5625

$ws.Range("H126").Value = 2075.75
$ws.Range("I126").Value = 4500
$ws.Range("K126").Value = 13500
$ws.Range("M126").Value = -11030

$ws.Range("H132").Value = 16135637
$ws.Range("I132").Value = 20002414
$ws.Range("K132").Value = 60007242
$ws.Range("M132").Value = -60004712

$ws.Range("H136").Value = 38463620
$ws.Range("I136").Value = 38463620
$ws.Range("K136").Value = 115390860
$ws.Range("M136").Value = -115388310
